$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) New header cells K2 / L2 (Eigenerosion / Eigensoftening)
# ------------------------------------------------------------------
$ws.Range("K2").Value = "Eigenerosion"
$ws.Range("L2").Value = "Eigensoftening"

# ------------------------------------------------------------------
# 2) New data rows 45-51 (Number / Material property / Definition).
#    Values are written in the same order the workbook's author typed
#    them so new shared-string entries land at the expected indices.
# ------------------------------------------------------------------
$ws.Range("B45").Value = 43
$ws.Range("D45").Value = "C EPSILON"
$ws.Range("C45").Value = "CEPS"

$ws.Range("B47").Value = 45
$ws.Range("C47").Value = "WC"
$ws.Range("D47").Value = "Wc"

$ws.Range("B48").Value = 46
$ws.Range("C48").Value = "FT"
$ws.Range("D48").Value = "Ft"

$ws.Range("B49").Value = 47
$ws.Range("C49").Value = "WC_P"
$ws.Range("D49").Value = "WC middle point"

$ws.Range("B50").Value = 48
$ws.Range("C50").Value = "FT_P"
$ws.Range("D50").Value = "FT middle point"

$ws.Range("B46").Value = 44
$ws.Range("C46").Value = "GC"
$ws.Range("D46").Value = "Gc"

$ws.Range("B51").Value = 49
$ws.Range("C51").Value = "D"
$ws.Range("D51").Value = "Aggregates Size"

# ------------------------------------------------------------------
# 3) "X" flag cells in the new K (Eigenerosion) / L (Eigensoftening)
#    columns for rows 45-51.
# ------------------------------------------------------------------
$ws.Range("K45").Value = "X"
$ws.Range("L45").Value = "X"
$ws.Range("K46").Value = "X"
$ws.Range("L47").Value = "X"
$ws.Range("L48").Value = "X"
$ws.Range("L49").Value = "X"
$ws.Range("L50").Value = "X"
$ws.Range("L51").Value = "X"

# ------------------------------------------------------------------
# 4) Formatting: mirror the alternating row styles already used by the
#    table (columns B/C/D) and the centered style used by the flag
#    columns, by copying formats from the existing rows/cells.
# ------------------------------------------------------------------
$ws.Range("B43").Copy()
$ws.Range("B45").PasteSpecial(-4122)
$ws.Range("B47").PasteSpecial(-4122)
$ws.Range("B49").PasteSpecial(-4122)
$ws.Range("B51").PasteSpecial(-4122)

$ws.Range("B44").Copy()
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("B48").PasteSpecial(-4122)
$ws.Range("B50").PasteSpecial(-4122)

$ws.Range("C43").Copy()
$ws.Range("C45:D51").PasteSpecial(-4122)

$ws.Range("I43").Copy()
$ws.Range("K45:L51").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 5) New column widths for K and L (closest achievable to the
#    author's saved best-fit widths).
# ------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 10.83
$ws.Columns.Item(12).ColumnWidth = 12.5

# ------------------------------------------------------------------
# 6) Restore the view: frozen header pane scrolled down plus the final
#    selection on L52.
# ------------------------------------------------------------------
$ws.Range("A34").Select()
$ws.Range("L52").Select()
